$wb = $excel.ActiveWorkbook

$loginsSheet = $wb.Worksheets.Item("Logins")
$productsSheet = $wb.Worksheets.Item("Products")

# ---- Logins sheet: add new rows 3-5 ----

# First, stamp the existing data style (font 1) onto the new range so that
# typed values pick up the same formatting as the rest of the table.
$loginsSheet.Range("A2").Copy()
$loginsSheet.Range("A3:C5").PasteSpecial(-4122)

# Row 3
$loginsSheet.Range("A3").Value = "locked"
$loginsSheet.Range("B3").Value = "secret"
$loginsSheet.Range("C3").Value = "Epic sadface: Username and password do not match any user in this service"

# Row 4 (blank-but-text cells, entered the way Excel represents a manually
# quote-prefixed empty entry: apostrophe only)
$loginsSheet.Range("A4").Value = "'"
$loginsSheet.Range("B4").Value = "'"
$loginsSheet.Range("C4").Value = "Epic sadface: Username is required"

# Row 5
$loginsSheet.Range("A5").Value = "locked_out_user"
$loginsSheet.Range("B5").Value = "'"
$loginsSheet.Range("C5").Value = "Epic sadface: Password is required"

# Column widths for Logins sheet
$loginsSheet.Columns.Item(1).ColumnWidth = 19.6667
$loginsSheet.Columns.Item(2).ColumnWidth = 14.1667
$loginsSheet.Columns.Item(3).ColumnWidth = 20.8333

# ---- Products sheet: apply currency number format to prices ----
$productsSheet.Range("B2:B7").NumberFormat = """$""#,##0.00"

# Column width for Products sheet
$productsSheet.Columns.Item(1).ColumnWidth = 20.5
